$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert J2, M2, P2, S2 from formulas to static values (drop the <f> element,
# keep the already-computed numeric <v> result unchanged).
$ws.Range("J2").Value2 = $ws.Range("J2").Value2
$ws.Range("M2").Value2 = $ws.Range("M2").Value2
$ws.Range("P2").Value2 = $ws.Range("P2").Value2
$ws.Range("S2").Value2 = $ws.Range("S2").Value2

# Clear out the now-unused totals row (J3, M3, P3, S3) while keeping their styles.
$ws.Range("J3").ClearContents()
$ws.Range("M3").ClearContents()
$ws.Range("P3").ClearContents()
$ws.Range("S3").ClearContents()

# Update the view: scroll to column E and move the active selection to M9.
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("M9").Select()
